# Generate Report for Handoff
#
# A new handoff cycle ran for "b.md": a fresh handoff package (.xlf) was
# generated for each target language. As a result:
#   - b.md's Status flips from "Handed back: in sync with en-US" to
#     "Ready for handoff" (Overview sheet + both language sheets).
#   - A new "Latest Handoff File" / "Latest Handoff Datetime" is recorded
#     on each language sheet.
#   - Because the previously handed-back translation is now stale relative
#     to the fresh handoff, "Content Duplicate" flips True -> False and an
#     "Error Detail" message is recorded on each language sheet.
#   - Column P ("Error Detail") is widened now that it holds real text.

$wb = $excel.ActiveWorkbook

# Helper: write literal text into a cell without letting Excel's automatic
# type inference turn look-alike tokens (e.g. "True"/"False") into a real
# Boolean. We stage the text as a formula string literal in a scratch cell,
# copy it, and Paste-Special "Values" into the destination -- this keeps the
# cell's type as a plain (shared) string, matching the source data.
function Set-LiteralText {
    param($range, [string]$text)

    $escaped = $text.Replace('"', '""')
    $scratch = $range.Worksheet.Range("ZZ1000")
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $scratch.Clear() | Out-Null
}

$newStatus   = "Ready for handoff"
$newDateTime = "2016-08-23 10:38:45"

# ---------------------------------------------------------------------
# Overview sheet: b.md row (row 3) -> Status columns (E, F) + date (G)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("G3").Value = $newDateTime

# ---------------------------------------------------------------------
# zh-cn sheet: b.md row (row 3)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(16).ColumnWidth = 40

$zhcn.Range("C3").Value = $newStatus
Set-LiteralText $zhcn.Range("F3") "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-23 10:38:40"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bfecaa25e74bf8730fcc63c878f88b807976b7b8/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a81648b97468f7f70f592fce9908cf58baaea6f5/e2e/b.md."

# ---------------------------------------------------------------------
# de-de sheet: b.md row (row 3)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(16).ColumnWidth = 40

$dede.Range("C3").Value = $newStatus
Set-LiteralText $dede.Range("F3") "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = $newDateTime
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bfecaa25e74bf8730fcc63c878f88b807976b7b8/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a81648b97468f7f70f592fce9908cf58baaea6f5/e2e/b.md."
